$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 6 new rows right after the header row (new rows 2-7), ---
# --- pushing the existing data rows down by 6 ---
$ws.Range("A2:C7").Insert()
$ws.Range("A2:C7").ClearFormats()

$topData = @(
    @(-0.3537254333496094, 0.2088937759399415, 0.5035260319709778),
    @(-0.4043011069297791, 0.3229363560676575, 0.4749223440885544),
    @(-0.2195036411285398, 0.2696369886398314, 0.4450621306896209),
    @(-0.5155707597732552, 0.2643678188323975, 0.5658968165516856),
    @(-0.4721715450286855, 0.2206716537475583, 0.4629700779914848),
    @(-0.2396689057350159, 0.2092438936233521, 0.5360905304551128)
)

for ($i = 0; $i -lt $topData.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $topData[$i][0]
    $ws.Cells.Item($r, 2).Value = $topData[$i][1]
    $ws.Cells.Item($r, 3).Value = $topData[$i][2]
}

# --- Append 4 new rows at the bottom (rows 28-31) ---
$bottomData = @(
    @(-0.1972274780273441, 0.3113194406032568, 0.5224930047988898),
    @(-0.2138409614562988, 0.2108606994152066, 0.46574055776),
    @(-0.09080266952514603, 0.2757070064544683, 0.4601370841264726),
    @(-0.04834830760955861, 0.3611972928047176, 0.5197352617979051)
)

for ($i = 0; $i -lt $bottomData.Length; $i++) {
    $r = 28 + $i
    $ws.Cells.Item($r, 1).Value = $bottomData[$i][0]
    $ws.Cells.Item($r, 2).Value = $bottomData[$i][1]
    $ws.Cells.Item($r, 3).Value = $bottomData[$i][2]
}
